$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row 6 for the nested pptx test entry
$ws.Range("B6").Value = "nested\test.pptx"
$ws.Range("C6").Value = "Van Jones"
$ws.Range("D6").Value = "Testing for nested pptx"

# Fix the typo in row 5 Title cell
$ws.Range("A5").Value = "Elder Scrolls Legends now on Ipad"

# Finish row 6 with the title text
$ws.Range("A6").Value = "Now 10% more awesome"

# Widen column B slightly to fit new content
$ws.Columns.Item(2).ColumnWidth = 16.140625

# Update the active selection to mirror the authored change
$ws.Range("E7").Select()
